# Insert a new data row at row 42 in the "Hortaliza, Feria Lagunitas de
# Puerto Montt - Ciboulette" sheet. This shifts the existing rows 42..130
# down to 43..131 (dimension grows from A1:R130 to A1:R131) and the new
# row 42 is populated with a fresh record for this market/category.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 42:130 down by one row, opening up a blank row 42.
$ws.Rows(42).Insert()

# Populate the newly inserted row 42 with the new record's data.
$ws.Cells.Item(42, 1).Value2  = 4
$ws.Cells.Item(42, 2).Value2  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(42, 3).Value2  = 'Los Lagos'
$ws.Cells.Item(42, 4).Value2  = 44498
$ws.Cells.Item(42, 5).Value2  = 10
$ws.Cells.Item(42, 6).Value2  = 100112039
$ws.Cells.Item(42, 7).Value2  = 'Ciboulette'
$ws.Cells.Item(42, 8).Value2  = 'Sin especificar'
$ws.Cells.Item(42, 9).Value2  = 'Primera'
$ws.Cells.Item(42, 10).Value2 = 240
$ws.Cells.Item(42, 11).Value2 = 2500
$ws.Cells.Item(42, 12).Value2 = 2500
$ws.Cells.Item(42, 13).Value2 = 2500
$ws.Cells.Item(42, 14).Value2 = '$/docena de atados'
$ws.Cells.Item(42, 15).Value2 = 'Región Metropolitana'
$ws.Cells.Item(42, 16).Value2 = 833
$ws.Cells.Item(42, 17).Value2 = 3
$ws.Cells.Item(42, 18).Value2 = 'Hortaliza'
